$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These "Price" cells hold numeric-looking text (e.g. "277.20", "0.06663")
# that must stay text so trailing/leading zeros survive, exactly as the
# source data had them. Temporarily mark them Text before typing the new
# value, then restore the Normal style so no extra formatting lingers.
$textCells = @(
    'D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D15', 'D16',
    'D18', 'D19', 'D20', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29',
    'D30', 'D31', 'D32', 'D33', 'D34', 'D37', 'D38', 'D39', 'D40', 'D41',
    'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Updated cryptos list (Price / Volume(1h) columns) ---

$ws.Range('D2').Value = '25.776.34'
$ws.Range('E2').Value = '  -3.98%  '
$ws.Range('D3').Value = '1.815.26'
$ws.Range('E3').Value = '  -3.12%  '
$ws.Range('D5').Value = '277.20'
$ws.Range('E5').Value = '  -7.92%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.5073'
$ws.Range('E7').Value = '  -4.79%  '
$ws.Range('D8').Value = '0.3524'
$ws.Range('E8').Value = '  -6.12%  '
$ws.Range('D9').Value = '44.58'
$ws.Range('E9').Value = '  -2.33%  '
$ws.Range('D10').Value = '0.06663'
$ws.Range('E10').Value = '  -7.17%  '
$ws.Range('D11').Value = '20.04'
$ws.Range('E11').Value = '  -7.35%  '
$ws.Range('D12').Value = '0.8275'
$ws.Range('E12').Value = '  -6.82%  '
$ws.Range('D13').Value = '0.07865'
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').Value = '1.798.59'
$ws.Range('E14').Value = '  -4.25%  '
$ws.Range('D15').Value = '5.075'
$ws.Range('E15').Value = '  -4.27%  '
$ws.Range('D16').Value = '87.48'
$ws.Range('E16').Value = '  -6.26%  '
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '14.12'
$ws.Range('E18').Value = '  -4.53%  '
$ws.Range('D19').Value = '0.000008044'
$ws.Range('E19').Value = '  -5.98%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').Value = '25.837.34'
$ws.Range('E21').Value = '  -4.04%  '
$ws.Range('D22').Value = '4.734'
$ws.Range('E22').Value = '  -5.10%  '
$ws.Range('E23').Value = '  -6.19%  '
$ws.Range('D24').Value = '6.094'
$ws.Range('E24').Value = '  -4.82%  '
$ws.Range('D25').Value = '142.15'
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('D26').Value = '2.193'
$ws.Range('E26').Value = '  -3.69%  '
$ws.Range('D27').Value = '1.673'
$ws.Range('E27').Value = '  -3.46%  '
$ws.Range('D28').Value = '17.10'
$ws.Range('E28').Value = '  -5.23%  '
$ws.Range('D29').Value = '109.38'
$ws.Range('E29').Value = '  -3.82%  '
$ws.Range('D30').Value = '4.340'
$ws.Range('E30').Value = '  -8.23%  '
$ws.Range('D31').Value = '4.236'
$ws.Range('E31').Value = '  -8.18%  '
$ws.Range('D32').Value = '0.08795'
$ws.Range('E32').Value = '  -3.72%  '
$ws.Range('D33').Value = '0.04880'
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('D34').Value = '0.7289'
$ws.Range('E34').Value = '  -10.66%  '
$ws.Range('E35').Value = '  -3.20%  '
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('D37').Value = '0.9999'
$ws.Range('D38').Value = '3.130'
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').Value = '2.377'
$ws.Range('E39').Value = '  -8.54%  '
$ws.Range('D40').Value = '0.01850'
$ws.Range('E40').Value = '  -5.20%  '
$ws.Range('D41').Value = '0.5176'
$ws.Range('E41').Value = '  -14.89%  '
$ws.Range('D42').Value = '0.9647'
$ws.Range('E42').Value = '  -9.92%  '
$ws.Range('D43').Value = '6.193'
$ws.Range('E43').Value = '  -6.09%  '
$ws.Range('D44').Value = '110.63'
$ws.Range('E44').Value = '  -3.98%  '
$ws.Range('D45').Value = '8.015'
$ws.Range('E45').Value = '  -10.00%  '
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = '0.4571'
$ws.Range('E47').Value = '  -11.35%  '
$ws.Range('D48').Value = '0.1363'
$ws.Range('E48').Value = '  -8.77%  '
$ws.Range('D49').Value = '36.58'
$ws.Range('E49').Value = '  -2.49%  '
$ws.Range('D50').Value = '9.289'
$ws.Range('E50').Value = '  -6.13%  '
$ws.Range('D51').Value = '1.498'
$ws.Range('E51').Value = '  -8.10%  '

# Restore default styling on the cells we temporarily forced to Text.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
